{"js": "// Apply the \"Forgotten Fable\" copy refresh:\n//  - New title/headline (appears twice: H1 heading + bold recap line)\n//  - \"What we like\" bullets reworded\n//  - \"What we don't like\" bullet reworded\n//  - Meta description (italic line) reworded\n\nconst body = context.document.body;\n\nconst replacements = [\n  {\n    find: \"Play Forgotten Fable for Free - Stunning Graphics & Exciting Features\",\n    replace: \"Play Forgotten Fable Free: Stunning Graphics and Big Wins\",\n  },\n  {\n    find: \"Stunning graphics\",\n    replace: \"Impeccable graphics and stunning design\",\n  },\n  {\n    find: \"Creative fantasy theme\",\n    replace: \"Unique fantasy theme featuring legendary villains\",\n  },\n  {\n    find: \"Generous winning Potential\",\n    replace: \"Special features with multipliers and free spins\",\n  },\n  {\n    find: \"Exciting special features\",\n    replace: \"High winning potential with up to 1,270 times the bet\",\n  },\n  {\n    find: \"Limited number of paylines\",\n    replace: \"Limited availability on some online casinos\",\n  },\n  {\n    find: \"Discover the Forgotten Fable slot game with its unique fantasy theme, stunning graphics, and exciting special features. Play for free and win big!\",\n    replace: \"Experience the unique fantasy theme and high winning potential of Forgotten Fable. Play for free now!\",\n  },\n];\n\nconst allResults = [];\nfor (const { find } of replacements) {\n  const results = body.search(find, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  allResults.push(results);\n}\nawait context.sync();\n\nfor (let i = 0; i < replacements.length; i++) {\n  const { replace } = replacements[i];\n  const results = allResults[i];\n  for (let j = 0; j < results.items.length; j++) {\n    results.items[j].insertText(replace, Word.InsertLocation.replace);\n  }\n}\nawait context.sync();\n", "ps1": "# Apply the \"Forgotten Fable\" copy refresh:\n#  - New title/headline (appears twice: H1 heading + bold recap line)\n#  - \"What we like\" bullets reworded\n#  - \"What we don't like\" bullet reworded\n#  - Meta description (italic line) reworded\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{\n        Find    = \"Play Forgotten Fable for Free - Stunning Graphics & Exciting Features\"\n        Replace = \"Play Forgotten Fable Free: Stunning Graphics and Big Wins\"\n    },\n    @{\n        Find    = \"Stunning graphics\"\n        Replace = \"Impeccable graphics and stunning design\"\n    },\n    @{\n        Find    = \"Creative fantasy theme\"\n        Replace = \"Unique fantasy theme featuring legendary villains\"\n    },\n    @{\n        Find    = \"Generous winning Potential\"\n        Replace = \"Special features with multipliers and free spins\"\n    },\n    @{\n        Find    = \"Exciting special features\"\n        Replace = \"High winning potential with up to 1,270 times the bet\"\n    },\n    @{\n        Find    = \"Limited number of paylines\"\n        Replace = \"Limited availability on some online casinos\"\n    },\n    @{\n        Find    = \"Discover the Forgotten Fable slot game with its unique fantasy theme, stunning graphics, and exciting special features. Play for free and win big!\"\n        Replace = \"Experience the unique fantasy theme and high winning potential of Forgotten Fable. Play for free now!\"\n    }\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    # Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards, MatchSoundsLike,\n    #         MatchAllWordForms, Forward, Wrap, Format, ReplaceWith, Replace)\n    # MatchCase=$true keeps this from e.g. matching the lowercase \"graphics\" inside\n    # unrelated sentences when the replacement target is capitalized, or vice versa.\n    $find.Execute($pair.Find, $true, $true, $false, $false, $false, $true, 1, $false, $pair.Replace, 2)\n}\n"}
